# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the refreshed counts published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheetId 1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 6352
$wsExhibit.Range("F5").Value  = 371
$wsExhibit.Range("F6").Value  = 54
$wsExhibit.Range("F9").Value  = 59
$wsExhibit.Range("F12").Value = 146
$wsExhibit.Range("F13").Value = 354
$wsExhibit.Range("F14").Value = 617
$wsExhibit.Range("F15").Value = 3090
$wsExhibit.Range("F16").Value = 8
$wsExhibit.Range("F18").Value = 1747

# --- Sheet "全部类型" (sheetId 4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6352
$wsAll.Range("F5").Value  = 371
$wsAll.Range("F6").Value  = 54
$wsAll.Range("F10").Value = 59
$wsAll.Range("F13").Value = 146
$wsAll.Range("F14").Value = 354
$wsAll.Range("F15").Value = 617
$wsAll.Range("F16").Value = 3090
$wsAll.Range("F17").Value = 8
$wsAll.Range("F19").Value = 1747
